$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 5411.8184  # H2: 4970.1665 -> 5411.8184
$ws.Cells.Item(2, 10).Value = 11433.223  # J2: 9374.817999999999 -> 11433.223
$ws.Cells.Item(2, 12).Value = 11433.223  # L2: 9374.817999999999 -> 11433.223
$ws.Cells.Item(2, 14).Value = -11659.223  # N2: -9600.817999999999 -> -11659.223
$ws.Cells.Item(3, 8).Value = 0  # H3: 28500 -> 0
$ws.Cells.Item(3, 10).Value = 0  # J3: 28500 -> 0
$ws.Cells.Item(3, 12).Value = 0  # L3: 28500 -> 0
$ws.Cells.Item(3, 14).Value = $null  # N3: -28728 -> (removed)
$ws.Cells.Item(4, 8).Value = 507.9375  # H4: 483.88235 -> 507.9375
$ws.Cells.Item(4, 9).Value = 597  # I4: 558.6923 -> 597
$ws.Cells.Item(4, 11).Value = 597  # K4: 558.6923 -> 597
$ws.Cells.Item(4, 13).Value = -483  # M4: -444.6923 -> -483
$ws.Cells.Item(19, 8).Value = 2462.5715  # H19: 2806.3333 -> 2462.5715
$ws.Cells.Item(19, 9).Value = 2462.5715  # I19: 2806.3333 -> 2462.5715
$ws.Cells.Item(19, 11).Value = 2462.5715  # K19: 2806.3333 -> 2462.5715
$ws.Cells.Item(19, 13).Value = -2287.5715  # M19: -2631.3333 -> -2287.5715
$ws.Cells.Item(29, 8).Value = 7037.75  # H29: 7204.6523 -> 7037.75
$ws.Cells.Item(29, 10).Value = 7270.647  # J29: 7525.125 -> 7270.647
$ws.Cells.Item(29, 12).Value = 21811.941  # L29: 22575.375 -> 21811.941
$ws.Cells.Item(29, 14).Value = -22373.941  # N29: -23137.375 -> -22373.941
$ws.Cells.Item(33, 8).Value = 300.52942  # H33: 332.66666 -> 300.52942
$ws.Cells.Item(33, 9).Value = 266.15384  # I33: 303.72726 -> 266.15384
$ws.Cells.Item(33, 11).Value = 266.15384  # K33: 303.72726 -> 266.15384
$ws.Cells.Item(33, 13).Value = -37.15384  # M33: -74.72726 -> -37.15384
$ws.Cells.Item(40, 8).Value = 2604.8076  # H40: 2515.7036 -> 2604.8076
$ws.Cells.Item(40, 10).Value = 2471.6667  # J40: 2244.4 -> 2471.6667
$ws.Cells.Item(40, 12).Value = 2471.6667  # L40: 2244.4 -> 2471.6667
$ws.Cells.Item(40, 14).Value = -2821.6667  # N40: -2594.4 -> -2821.6667
$ws.Cells.Item(43, 8).Value = 33325  # H43: 27766.666 -> 33325
$ws.Cells.Item(43, 9).Value = 55000  # I43: 60000 -> 55000
$ws.Cells.Item(43, 11).Value = 55000  # K43: 60000 -> 55000
$ws.Cells.Item(43, 13).Value = -54931  # M43: -59931 -> -54931
$ws.Cells.Item(63, 8).Value = 80000  # H63: 0 -> 80000
$ws.Cells.Item(63, 10).Value = 80000  # J63: 0 -> 80000
$ws.Cells.Item(63, 12).Value = 80000  # L63: 0 -> 80000
$ws.Cells.Item(63, 14).Value = -81248  # N63: None -> -81248
$ws.Cells.Item(66, 8).Value = 80000  # H66: 0 -> 80000
$ws.Cells.Item(66, 10).Value = 80000  # J66: 0 -> 80000
$ws.Cells.Item(66, 12).Value = 240000  # L66: 0 -> 240000
$ws.Cells.Item(66, 14).Value = -246240  # N66: None -> -246240
$ws.Cells.Item(80, 8).Value = 2312.5908  # H80: 2507.6843 -> 2312.5908
$ws.Cells.Item(80, 9).Value = 1805.5625  # I80: 1876.5333 -> 1805.5625
$ws.Cells.Item(80, 10).Value = 3664.6667  # J80: 4874.5 -> 3664.6667
$ws.Cells.Item(80, 11).Value = 5416.6875  # K80: 5629.5999 -> 5416.6875
$ws.Cells.Item(80, 12).Value = 10994.0001  # L80: 14623.5 -> 10994.0001
$ws.Cells.Item(80, 13).Value = -4418.6875  # M80: -4631.5999 -> -4418.6875
$ws.Cells.Item(80, 14).Value = -12990.0001  # N80: -16619.5 -> -12990.0001
$ws.Cells.Item(83, 8).Value = 2312.5908  # H83: 2507.6843 -> 2312.5908
$ws.Cells.Item(83, 9).Value = 1805.5625  # I83: 1876.5333 -> 1805.5625
$ws.Cells.Item(83, 10).Value = 3664.6667  # J83: 4874.5 -> 3664.6667
$ws.Cells.Item(83, 11).Value = 16250.0625  # K83: 16888.7997 -> 16250.0625
$ws.Cells.Item(83, 12).Value = 32982.0003  # L83: 43870.5 -> 32982.0003
$ws.Cells.Item(83, 13).Value = -11258.0625  # M83: -11896.7997 -> -11258.0625
$ws.Cells.Item(83, 14).Value = -42966.0003  # N83: -53854.5 -> -42966.0003
$ws.Cells.Item(86, 8).Value = 2446.5557  # H86: 2517.1428 -> 2446.5557
$ws.Cells.Item(86, 9).Value = 2467  # I86: 2574 -> 2467
$ws.Cells.Item(86, 11).Value = 2467  # K86: 2574 -> 2467
$ws.Cells.Item(86, 13).Value = -1344  # M86: -1451 -> -1344
$ws.Cells.Item(89, 8).Value = 2446.5557  # H89: 2517.1428 -> 2446.5557
$ws.Cells.Item(89, 9).Value = 2467  # I89: 2574 -> 2467
$ws.Cells.Item(89, 11).Value = 12335  # K89: 12870 -> 12335
$ws.Cells.Item(89, 13).Value = -6719  # M89: -7254 -> -6719
$ws.Cells.Item(92, 8).Value = 4181.4  # H92: 3801.1667 -> 4181.4
$ws.Cells.Item(92, 9).Value = 4181.4  # I92: 3801.1667 -> 4181.4
$ws.Cells.Item(92, 11).Value = 4181.4  # K92: 3801.1667 -> 4181.4
$ws.Cells.Item(92, 13).Value = -2933.4  # M92: -2553.1667 -> -2933.4
$ws.Cells.Item(96, 8).Value = 12531.286  # H96: 6514.7856 -> 12531.286
$ws.Cells.Item(96, 9).Value = 585  # I96: 332.25 -> 585
$ws.Cells.Item(96, 10).Value = 17309.8  # J96: 14758.167 -> 17309.8
$ws.Cells.Item(96, 11).Value = 1755  # K96: 996.75 -> 1755
$ws.Cells.Item(96, 12).Value = 51929.39999999999  # L96: 44274.501 -> 51929.39999999999
$ws.Cells.Item(96, 13).Value = -382  # M96: 376.25 -> -382
$ws.Cells.Item(96, 14).Value = -54675.39999999999  # N96: -47020.501 -> -54675.39999999999
$ws.Cells.Item(97, 8).Value = 27149.5  # H97: 13598.444 -> 27149.5
$ws.Cells.Item(97, 10).Value = 27149.5  # J97: 13598.444 -> 27149.5
$ws.Cells.Item(97, 12).Value = 81448.5  # L97: 40795.33199999999 -> 81448.5
$ws.Cells.Item(97, 14).Value = -82440.5  # N97: -41787.33199999999 -> -82440.5
$ws.Cells.Item(99, 8).Value = 404.7143  # H99: 441.7143 -> 404.7143
$ws.Cells.Item(99, 9).Value = 408.6  # I99: 415.33334 -> 408.6
$ws.Cells.Item(99, 10).Value = 395  # J99: 600 -> 395
$ws.Cells.Item(99, 11).Value = 1225.8  # K99: 1246.00002 -> 1225.8
$ws.Cells.Item(99, 12).Value = 1185  # L99: 1800 -> 1185
$ws.Cells.Item(99, 13).Value = 272.1999999999998  # M99: 251.9999800000001 -> 272.1999999999998
$ws.Cells.Item(99, 14).Value = -4181  # N99: -4796 -> -4181
$ws.Cells.Item(101, 8).Value = 1722  # H101: 2051.2222 -> 1722
$ws.Cells.Item(101, 9).Value = 1656  # I101: 2013.1428 -> 1656
$ws.Cells.Item(101, 10).Value = 2184  # J101: 2184.5 -> 2184
$ws.Cells.Item(101, 11).Value = 4968  # K101: 6039.428400000001 -> 4968
$ws.Cells.Item(101, 12).Value = 6552  # L101: 6553.5 -> 6552
$ws.Cells.Item(101, 13).Value = -3346  # M101: -4417.428400000001 -> -3346
$ws.Cells.Item(101, 14).Value = -9796  # N101: -9797.5 -> -9796
$ws.Cells.Item(102, 8).Value = 0  # H102: 28500 -> 0
$ws.Cells.Item(102, 10).Value = 0  # J102: 28500 -> 0
$ws.Cells.Item(102, 12).Value = 0  # L102: 28500 -> 0
$ws.Cells.Item(102, 14).Value = $null  # N102: -34990 -> (removed)
$ws.Cells.Item(106, 8).Value = 23511.928  # H106: 21010.125 -> 23511.928
$ws.Cells.Item(106, 9).Value = 27281.666  # I106: 23883.928 -> 27281.666
$ws.Cells.Item(106, 11).Value = 27281.666  # K106: 23883.928 -> 27281.666
$ws.Cells.Item(106, 13).Value = -26650.666  # M106: -23252.928 -> -26650.666
$ws.Cells.Item(131, 8).Value = 3366  # H131: 3186.4443 -> 3366
$ws.Cells.Item(131, 10).Value = 5000  # J131: 3375 -> 5000
$ws.Cells.Item(131, 12).Value = 15000  # L131: 10125 -> 15000
$ws.Cells.Item(131, 14).Value = -25080  # N131: -20205 -> -25080
$ws.Cells.Item(138, 8).Value = 2671.5134  # H138: 2617.8157 -> 2671.5134
$ws.Cells.Item(138, 9).Value = 1627.1875  # I138: 1627.25 -> 1627.1875
$ws.Cells.Item(138, 10).Value = 3467.1904  # J138: 3338.2273 -> 3467.1904
$ws.Cells.Item(138, 11).Value = 4881.5625  # K138: 4881.75 -> 4881.5625
$ws.Cells.Item(138, 12).Value = 10401.5712  # L138: 10014.6819 -> 10401.5712
$ws.Cells.Item(138, 13).Value = 258.4375  # M138: 258.25 -> 258.4375
$ws.Cells.Item(138, 14).Value = -20681.5712  # N138: -20294.6819 -> -20681.5712

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 26881.465  # H32: 27840.264 -> 26881.465
$ws.Cells.Item(32, 9).Value = 27473.547  # I32: 28565.676 -> 27473.547
$ws.Cells.Item(32, 10).Value = 2014  # J32: 1000 -> 2014
$ws.Cells.Item(32, 11).Value = 27473.547  # K32: 28565.676 -> 27473.547
$ws.Cells.Item(32, 12).Value = 2014  # L32: 1000 -> 2014
$ws.Cells.Item(32, 13).Value = -27186.547  # M32: -28278.676 -> -27186.547
$ws.Cells.Item(32, 14).Value = -2588  # N32: -1574 -> -2588
$ws.Cells.Item(45, 8).Value = 9846.615  # H45: 14873.5 -> 9846.615
$ws.Cells.Item(45, 9).Value = 10250.5  # I45: 15855.429 -> 10250.5
$ws.Cells.Item(45, 10).Value = 5000  # J45: 8000 -> 5000
$ws.Cells.Item(45, 11).Value = 10250.5  # K45: 15855.429 -> 10250.5
$ws.Cells.Item(45, 12).Value = 5000  # L45: 8000 -> 5000
$ws.Cells.Item(45, 13).Value = -9873.5  # M45: -15478.429 -> -9873.5
$ws.Cells.Item(45, 14).Value = -5754  # N45: -8754 -> -5754
$ws.Cells.Item(61, 8).Value = 6767.8  # H61: 6774.65 -> 6767.8
$ws.Cells.Item(61, 9).Value = 3937.1333  # I61: 3946.2666 -> 3937.1333
$ws.Cells.Item(61, 11).Value = 3937.1333  # K61: 3946.2666 -> 3937.1333
$ws.Cells.Item(61, 13).Value = -3725.1333  # M61: -3734.2666 -> -3725.1333
$ws.Cells.Item(74, 8).Value = 328119.72  # H74: 317882.7 -> 328119.72
$ws.Cells.Item(74, 9).Value = 503510.6  # I74: 457818 -> 503510.6
$ws.Cells.Item(74, 10).Value = 9227.182000000001  # J74: 10025 -> 9227.182000000001
$ws.Cells.Item(74, 11).Value = 503510.6  # K74: 457818 -> 503510.6
$ws.Cells.Item(74, 12).Value = 9227.182000000001  # L74: 10025 -> 9227.182000000001
$ws.Cells.Item(74, 13).Value = -502636.6  # M74: -456944 -> -502636.6
$ws.Cells.Item(74, 14).Value = -10975.182  # N74: -11773 -> -10975.182
$ws.Cells.Item(77, 8).Value = 328119.72  # H77: 317882.7 -> 328119.72
$ws.Cells.Item(77, 9).Value = 503510.6  # I77: 457818 -> 503510.6
$ws.Cells.Item(77, 10).Value = 9227.182000000001  # J77: 10025 -> 9227.182000000001
$ws.Cells.Item(77, 11).Value = 2517553  # K77: 2289090 -> 2517553
$ws.Cells.Item(77, 12).Value = 46135.91  # L77: 50125 -> 46135.91
$ws.Cells.Item(77, 13).Value = -2513185  # M77: -2284722 -> -2513185
$ws.Cells.Item(77, 14).Value = -54871.91  # N77: -58861 -> -54871.91
$ws.Cells.Item(97, 8).Value = 2061017.5  # H97: 1519997.8 -> 2061017.5
$ws.Cells.Item(97, 9).Value = 2473104.2  # I97: 1952499.2 -> 2473104.2
$ws.Cells.Item(97, 10).Value = 583.3333  # J97: 150409.83 -> 583.3333
$ws.Cells.Item(97, 11).Value = 2473104.2  # K97: 1952499.2 -> 2473104.2
$ws.Cells.Item(97, 12).Value = 583.3333  # L97: 150409.83 -> 583.3333
$ws.Cells.Item(97, 13).Value = -2472608.2  # M97: -1952003.2 -> -2472608.2
$ws.Cells.Item(97, 14).Value = -1575.3333  # N97: -151401.83 -> -1575.3333
$ws.Cells.Item(132, 8).Value = 11498.333  # H132: 2711.1714 -> 11498.333
$ws.Cells.Item(132, 9).Value = 6000  # I132: 1824.4 -> 6000
$ws.Cells.Item(132, 10).Value = 14247.5  # J132: 8031.8 -> 14247.5
$ws.Cells.Item(132, 11).Value = 18000  # K132: 5473.200000000001 -> 18000
$ws.Cells.Item(132, 12).Value = 42742.5  # L132: 24095.4 -> 42742.5
$ws.Cells.Item(132, 13).Value = -15470  # M132: -2943.200000000001 -> -15470
$ws.Cells.Item(132, 14).Value = -47802.5  # N132: -29155.4 -> -47802.5
$ws.Cells.Item(136, 8).Value = 6767.8  # H136: 6774.65 -> 6767.8
$ws.Cells.Item(136, 9).Value = 3937.1333  # I136: 3946.2666 -> 3937.1333
$ws.Cells.Item(136, 11).Value = 11811.3999  # K136: 11838.7998 -> 11811.3999
$ws.Cells.Item(136, 13).Value = -9261.3999  # M136: -9288.799800000001 -> -9261.3999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(54, 8).Value = 6145.3335  # H54: 9399.666999999999 -> 6145.3335
$ws.Cells.Item(54, 9).Value = 960  # I54: 1599.5 -> 960
$ws.Cells.Item(54, 10).Value = 11330.667  # J54: 25000 -> 11330.667
$ws.Cells.Item(54, 11).Value = 960  # K54: 1599.5 -> 960
$ws.Cells.Item(54, 12).Value = 11330.667  # L54: 25000 -> 11330.667
$ws.Cells.Item(54, 13).Value = -476  # M54: -1115.5 -> -476
$ws.Cells.Item(54, 14).Value = -12298.667  # N54: -25968 -> -12298.667
$ws.Cells.Item(88, 8).Value = 32356.625  # H88: 32432 -> 32356.625
$ws.Cells.Item(88, 10).Value = 32356.625  # J88: 32432 -> 32356.625
$ws.Cells.Item(88, 12).Value = 32356.625  # L88: 32432 -> 32356.625
$ws.Cells.Item(88, 14).Value = -33168.625  # N88: -33244 -> -33168.625
$ws.Cells.Item(91, 8).Value = 32356.625  # H91: 32432 -> 32356.625
$ws.Cells.Item(91, 10).Value = 32356.625  # J91: 32432 -> 32356.625
$ws.Cells.Item(91, 12).Value = 32356.625  # L91: 32432 -> 32356.625
$ws.Cells.Item(91, 14).Value = -35164.625  # N91: -35240 -> -35164.625
$ws.Cells.Item(106, 8).Value = 31060  # H106: 45000 -> 31060
$ws.Cells.Item(106, 10).Value = 31060  # J106: 45000 -> 31060
$ws.Cells.Item(106, 12).Value = 31060  # L106: 45000 -> 31060
$ws.Cells.Item(106, 14).Value = -33584  # N106: -47524 -> -33584
$ws.Cells.Item(107, 8).Value = 2479.0952  # H107: 2552.7273 -> 2479.0952
$ws.Cells.Item(107, 9).Value = 2479.0952  # I107: 2552.7273 -> 2479.0952
$ws.Cells.Item(107, 11).Value = 2479.0952  # K107: 2552.7273 -> 2479.0952
$ws.Cells.Item(107, 13).Value = -559.0952000000002  # M107: -632.7273 -> -559.0952000000002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(15, 8).Value = 4966.3335  # H15: 6997 -> 4966.3335
$ws.Cells.Item(15, 10).Value = 950  # J15: 995 -> 950
$ws.Cells.Item(15, 12).Value = 950  # L15: 995 -> 950
$ws.Cells.Item(15, 14).Value = -1290  # N15: -1335 -> -1290
$ws.Cells.Item(28, 8).Value = 33737.5  # H28: 34750 -> 33737.5
$ws.Cells.Item(28, 10).Value = 33737.5  # J28: 34750 -> 33737.5
$ws.Cells.Item(28, 12).Value = 33737.5  # L28: 34750 -> 33737.5
$ws.Cells.Item(28, 14).Value = -34227.5  # N28: -35240 -> -34227.5
$ws.Cells.Item(31, 8).Value = 33336346  # H31: 38464852 -> 33336346
$ws.Cells.Item(31, 9).Value = 47620316  # I31: 55556916 -> 47620316
$ws.Cells.Item(31, 10).Value = 7077.222  # J31: 7710.5 -> 7077.222
$ws.Cells.Item(31, 11).Value = 47620316  # K31: 55556916 -> 47620316
$ws.Cells.Item(31, 12).Value = 7077.222  # L31: 7710.5 -> 7077.222
$ws.Cells.Item(31, 13).Value = -47620021  # M31: -55556621 -> -47620021
$ws.Cells.Item(31, 14).Value = -7667.222  # N31: -8300.5 -> -7667.222
$ws.Cells.Item(34, 8).Value = 33336346  # H34: 38464852 -> 33336346
$ws.Cells.Item(34, 9).Value = 47620316  # I34: 55556916 -> 47620316
$ws.Cells.Item(34, 10).Value = 7077.222  # J34: 7710.5 -> 7077.222
$ws.Cells.Item(34, 11).Value = 47620316  # K34: 55556916 -> 47620316
$ws.Cells.Item(34, 12).Value = 7077.222  # L34: 7710.5 -> 7077.222
$ws.Cells.Item(34, 13).Value = -47620114  # M34: -55556714 -> -47620114
$ws.Cells.Item(34, 14).Value = -7481.222  # N34: -8114.5 -> -7481.222
$ws.Cells.Item(43, 8).Value = 87052.37  # H43: 98478.05 -> 87052.37
$ws.Cells.Item(43, 10).Value = 87052.37  # J43: 98478.05 -> 87052.37
$ws.Cells.Item(43, 12).Value = 87052.37  # L43: 98478.05 -> 87052.37
$ws.Cells.Item(43, 14).Value = -87420.37  # N43: -98846.05 -> -87420.37
$ws.Cells.Item(58, 8).Value = 5856.9414  # H58: 5865.0586 -> 5856.9414
$ws.Cells.Item(58, 9).Value = 3620.6155  # I58: 3631.2307 -> 3620.6155
$ws.Cells.Item(58, 11).Value = 3620.6155  # K58: 3631.2307 -> 3620.6155
$ws.Cells.Item(58, 13).Value = -3417.6155  # M58: -3428.2307 -> -3417.6155
$ws.Cells.Item(99, 8).Value = 3234.2  # H99: 3350 -> 3234.2
$ws.Cells.Item(99, 9).Value = 3234.2  # I99: 3437.5 -> 3234.2
$ws.Cells.Item(99, 10).Value = 0  # J99: 3000 -> 0
$ws.Cells.Item(99, 11).Value = 3234.2  # K99: 3437.5 -> 3234.2
$ws.Cells.Item(99, 12).Value = 0  # L99: 3000 -> 0
$ws.Cells.Item(99, 13).Value = -1736.2  # M99: -1939.5 -> -1736.2
$ws.Cells.Item(99, 14).Value = $null  # N99: -5996 -> (removed)
$ws.Cells.Item(101, 8).Value = 87052.37  # H101: 98478.05 -> 87052.37
$ws.Cells.Item(101, 10).Value = 87052.37  # J101: 98478.05 -> 87052.37
$ws.Cells.Item(101, 12).Value = 87052.37  # L101: 98478.05 -> 87052.37
$ws.Cells.Item(101, 14).Value = -93542.37  # N101: -104968.05 -> -93542.37
$ws.Cells.Item(105, 8).Value = 1720.1923  # H105: 1872.95 -> 1720.1923
$ws.Cells.Item(105, 9).Value = 1104.1111  # I105: 1192.1538 -> 1104.1111
$ws.Cells.Item(105, 10).Value = 3106.375  # J105: 3137.2856 -> 3106.375
$ws.Cells.Item(105, 11).Value = 1104.1111  # K105: 1192.1538 -> 1104.1111
$ws.Cells.Item(105, 12).Value = 3106.375  # L105: 3137.2856 -> 3106.375
$ws.Cells.Item(105, 13).Value = 642.8888999999999  # M105: 554.8462 -> 642.8888999999999
$ws.Cells.Item(105, 14).Value = -6600.375  # N105: -6631.2856 -> -6600.375
$ws.Cells.Item(126, 8).Value = 3234.2  # H126: 3350 -> 3234.2
$ws.Cells.Item(126, 9).Value = 3234.2  # I126: 3437.5 -> 3234.2
$ws.Cells.Item(126, 10).Value = 0  # J126: 3000 -> 0
$ws.Cells.Item(126, 11).Value = 9702.599999999999  # K126: 10312.5 -> 9702.599999999999
$ws.Cells.Item(126, 12).Value = 0  # L126: 9000 -> 0
$ws.Cells.Item(126, 13).Value = -7232.599999999999  # M126: -7842.5 -> -7232.599999999999
$ws.Cells.Item(126, 14).Value = $null  # N126: -13940 -> (removed)
$ws.Cells.Item(136, 8).Value = 5856.9414  # H136: 5865.0586 -> 5856.9414
$ws.Cells.Item(136, 9).Value = 3620.6155  # I136: 3631.2307 -> 3620.6155
$ws.Cells.Item(136, 11).Value = 10861.8465  # K136: 10893.6921 -> 10861.8465
$ws.Cells.Item(136, 13).Value = -8311.8465  # M136: -8343.6921 -> -8311.8465

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 786.9  # H5: 912.55554 -> 786.9
$ws.Cells.Item(5, 9).Value = 566.2857  # I5: 661.6 -> 566.2857
$ws.Cells.Item(5, 10).Value = 1301.6666  # J5: 1226.25 -> 1301.6666
$ws.Cells.Item(5, 11).Value = 1698.8571  # K5: 1984.8 -> 1698.8571
$ws.Cells.Item(5, 12).Value = 3904.9998  # L5: 3678.75 -> 3904.9998
$ws.Cells.Item(5, 13).Value = -1586.8571  # M5: -1872.8 -> -1586.8571
$ws.Cells.Item(5, 14).Value = -4128.9998  # N5: -3902.75 -> -4128.9998
$ws.Cells.Item(7, 8).Value = 71.5625  # H7: 74.166664 -> 71.5625
$ws.Cells.Item(7, 10).Value = 94.14286  # J7: 94.333336 -> 94.14286
$ws.Cells.Item(7, 12).Value = 282.42858  # L7: 283.000008 -> 282.42858
$ws.Cells.Item(7, 14).Value = -506.42858  # N7: -507.000008 -> -506.42858
$ws.Cells.Item(11, 8).Value = 6659.4443  # H11: 6309 -> 6659.4443
$ws.Cells.Item(11, 9).Value = 7173.1875  # I11: 6751.294 -> 7173.1875
$ws.Cells.Item(11, 11).Value = 21519.5625  # K11: 20253.882 -> 21519.5625
$ws.Cells.Item(11, 13).Value = -21379.5625  # M11: -20113.882 -> -21379.5625
$ws.Cells.Item(68, 8).Value = 3015  # H68: 2670.7144 -> 3015
$ws.Cells.Item(68, 9).Value = 772.5  # I68: 773.75 -> 772.5
$ws.Cells.Item(68, 10).Value = 7500  # J68: 5200 -> 7500
$ws.Cells.Item(68, 11).Value = 2317.5  # K68: 2321.25 -> 2317.5
$ws.Cells.Item(68, 12).Value = 22500  # L68: 15600 -> 22500
$ws.Cells.Item(68, 13).Value = -1506.5  # M68: -1510.25 -> -1506.5
$ws.Cells.Item(68, 14).Value = -24122  # N68: -17222 -> -24122
$ws.Cells.Item(71, 8).Value = 3015  # H71: 2670.7144 -> 3015
$ws.Cells.Item(71, 9).Value = 772.5  # I71: 773.75 -> 772.5
$ws.Cells.Item(71, 10).Value = 7500  # J71: 5200 -> 7500
$ws.Cells.Item(71, 11).Value = 6952.5  # K71: 6963.75 -> 6952.5
$ws.Cells.Item(71, 12).Value = 67500  # L71: 46800 -> 67500
$ws.Cells.Item(71, 13).Value = -2896.5  # M71: -2907.75 -> -2896.5
$ws.Cells.Item(71, 14).Value = -75612  # N71: -54912 -> -75612
$ws.Cells.Item(86, 8).Value = 183.33333  # H86: 304 -> 183.33333
$ws.Cells.Item(86, 9).Value = 150  # I86: 225 -> 150
$ws.Cells.Item(86, 10).Value = 250  # J86: 462 -> 250
$ws.Cells.Item(86, 11).Value = 450  # K86: 675 -> 450
$ws.Cells.Item(86, 12).Value = 750  # L86: 1386 -> 750
$ws.Cells.Item(86, 13).Value = 736  # M86: 511 -> 736
$ws.Cells.Item(86, 14).Value = -3122  # N86: -3758 -> -3122
$ws.Cells.Item(89, 8).Value = 183.33333  # H89: 304 -> 183.33333
$ws.Cells.Item(89, 9).Value = 150  # I89: 225 -> 150
$ws.Cells.Item(89, 10).Value = 250  # J89: 462 -> 250
$ws.Cells.Item(89, 11).Value = 1350  # K89: 2025 -> 1350
$ws.Cells.Item(89, 12).Value = 2250  # L89: 4158 -> 2250
$ws.Cells.Item(89, 13).Value = 4578  # M89: 3903 -> 4578
$ws.Cells.Item(89, 14).Value = -14106  # N89: -16014 -> -14106
$ws.Cells.Item(116, 8).Value = 2109.6  # H116: 2262.25 -> 2109.6
$ws.Cells.Item(116, 9).Value = 2109.6  # I116: 2262.25 -> 2109.6
$ws.Cells.Item(116, 11).Value = 6328.799999999999  # K116: 6786.75 -> 6328.799999999999
$ws.Cells.Item(116, 13).Value = -2886.799999999999  # M116: -3344.75 -> -2886.799999999999
$ws.Cells.Item(117, 8).Value = 1645.25  # H117: 2862.5 -> 1645.25
$ws.Cells.Item(117, 9).Value = 250.66667  # I117: 209 -> 250.66667
$ws.Cells.Item(117, 10).Value = 2482  # J117: 5516 -> 2482
$ws.Cells.Item(117, 11).Value = 752.00001  # K117: 627 -> 752.00001
$ws.Cells.Item(117, 12).Value = 7446  # L117: 16548 -> 7446
$ws.Cells.Item(117, 13).Value = 2689.99999  # M117: 2815 -> 2689.99999
$ws.Cells.Item(117, 14).Value = -14330  # N117: -23432 -> -14330
$ws.Cells.Item(127, 8).Value = 9500  # H127: 11325 -> 9500
$ws.Cells.Item(127, 10).Value = 9500  # J127: 11325 -> 9500
$ws.Cells.Item(127, 12).Value = 28500  # L127: 33975 -> 28500
$ws.Cells.Item(127, 14).Value = -38420  # N127: -43895 -> -38420
$ws.Cells.Item(131, 8).Value = 15880338  # H131: 17551848 -> 15880338
$ws.Cells.Item(131, 9).Value = 83334260  # I131: 111112010 -> 83334260
$ws.Cells.Item(131, 10).Value = 8827.471  # J131: 9316.6875 -> 8827.471
$ws.Cells.Item(131, 11).Value = 250002780  # K131: 333336030 -> 250002780
$ws.Cells.Item(131, 12).Value = 26482.413  # L131: 27950.0625 -> 26482.413
$ws.Cells.Item(131, 13).Value = -249997740  # M131: -333330990 -> -249997740
$ws.Cells.Item(131, 14).Value = -36562.413  # N131: -38030.0625 -> -36562.413
$ws.Cells.Item(135, 8).Value = 786.9  # H135: 912.55554 -> 786.9
$ws.Cells.Item(135, 9).Value = 566.2857  # I135: 661.6 -> 566.2857
$ws.Cells.Item(135, 10).Value = 1301.6666  # J135: 1226.25 -> 1301.6666
$ws.Cells.Item(135, 11).Value = 5096.571300000001  # K135: 5954.400000000001 -> 5096.571300000001
$ws.Cells.Item(135, 12).Value = 11714.9994  # L135: 11036.25 -> 11714.9994
$ws.Cells.Item(135, 13).Value = -2561.571300000001  # M135: -3419.400000000001 -> -2561.571300000001
$ws.Cells.Item(135, 14).Value = -16784.9994  # N135: -16106.25 -> -16784.9994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 20799.875  # H70: 19043.777 -> 20799.875
$ws.Cells.Item(70, 9).Value = 16842.715  # I70: 15361.75 -> 16842.715
$ws.Cells.Item(70, 11).Value = 16842.715  # K70: 15361.75 -> 16842.715
$ws.Cells.Item(70, 13).Value = -16572.715  # M70: -15091.75 -> -16572.715
$ws.Cells.Item(73, 8).Value = 20799.875  # H73: 19043.777 -> 20799.875
$ws.Cells.Item(73, 9).Value = 16842.715  # I73: 15361.75 -> 16842.715
$ws.Cells.Item(73, 11).Value = 16842.715  # K73: 15361.75 -> 16842.715
$ws.Cells.Item(73, 13).Value = -15906.715  # M73: -14425.75 -> -15906.715
$ws.Cells.Item(98, 8).Value = 29497  # H98: 35329.668 -> 29497
$ws.Cells.Item(98, 10).Value = 29497  # J98: 35329.668 -> 29497
$ws.Cells.Item(98, 12).Value = 29497  # L98: 35329.668 -> 29497
$ws.Cells.Item(98, 14).Value = -35487  # N98: -41319.668 -> -35487
$ws.Cells.Item(101, 8).Value = 35171.332  # H101: 33328.5 -> 35171.332
$ws.Cells.Item(101, 10).Value = 35171.332  # J101: 33328.5 -> 35171.332
$ws.Cells.Item(101, 12).Value = 35171.332  # L101: 33328.5 -> 35171.332
$ws.Cells.Item(101, 14).Value = -41661.332  # N101: -39818.5 -> -41661.332
$ws.Cells.Item(107, 8).Value = 823.5454999999999  # H107: 873.5 -> 823.5454999999999
$ws.Cells.Item(107, 9).Value = 862.2222  # I107: 904.5 -> 862.2222
$ws.Cells.Item(107, 10).Value = 649.5  # J107: 749.5 -> 649.5
$ws.Cells.Item(107, 11).Value = 862.2222  # K107: 904.5 -> 862.2222
$ws.Cells.Item(107, 12).Value = 649.5  # L107: 749.5 -> 649.5
$ws.Cells.Item(107, 13).Value = 1057.7778  # M107: 1015.5 -> 1057.7778
$ws.Cells.Item(107, 14).Value = -4489.5  # N107: -4589.5 -> -4489.5
$ws.Cells.Item(122, 8).Value = 7333.469  # H122: 7469.729 -> 7333.469
$ws.Cells.Item(122, 9).Value = 8581.281999999999  # I122: 9001.919 -> 8581.281999999999
$ws.Cells.Item(122, 10).Value = 2467  # J122: 2316 -> 2467
$ws.Cells.Item(122, 11).Value = 25743.846  # K122: 27005.757 -> 25743.846
$ws.Cells.Item(122, 12).Value = 7401  # L122: 6948 -> 7401
$ws.Cells.Item(122, 13).Value = -23293.846  # M122: -24555.757 -> -23293.846
$ws.Cells.Item(122, 14).Value = -12301  # N122: -11848 -> -12301
$ws.Cells.Item(132, 8).Value = 6953.25  # H132: 6054.933 -> 6953.25
$ws.Cells.Item(132, 9).Value = 3920.4285  # I132: 3482.8 -> 3920.4285
$ws.Cells.Item(132, 11).Value = 11761.2855  # K132: 10448.4 -> 11761.2855
$ws.Cells.Item(132, 13).Value = -9231.2855  # M132: -7918.400000000001 -> -9231.2855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 9843.111000000001  # H132: 9444.223 -> 9843.111000000001
$ws.Cells.Item(132, 9).Value = 5463.3335  # I132: 4266.6665 -> 5463.3335
$ws.Cells.Item(132, 11).Value = 16390.0005  # K132: 12799.9995 -> 16390.0005
$ws.Cells.Item(132, 13).Value = -13860.0005  # M132: -10269.9995 -> -13860.0005

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6094.684  # H132: 4546.8623 -> 6094.684
$ws.Cells.Item(132, 9).Value = 3579.8  # I132: 2233.182 -> 3579.8
$ws.Cells.Item(132, 10).Value = 6992.857  # J132: 5960.778 -> 6992.857
$ws.Cells.Item(132, 11).Value = 10739.4  # K132: 6699.545999999999 -> 10739.4
$ws.Cells.Item(132, 12).Value = 20978.571  # L132: 17882.334 -> 20978.571
$ws.Cells.Item(132, 13).Value = -8209.400000000001  # M132: -4169.545999999999 -> -8209.400000000001
$ws.Cells.Item(132, 14).Value = -26038.571  # N132: -22942.334 -> -26038.571
$ws.Cells.Item(136, 8).Value = 3413.9333  # H136: 2756.6843 -> 3413.9333
$ws.Cells.Item(136, 9).Value = 1155.7273  # I136: 925.4 -> 1155.7273
$ws.Cells.Item(136, 11).Value = 3467.1819  # K136: 2776.2 -> 3467.1819
$ws.Cells.Item(136, 13).Value = -917.1819  # M136: -226.1999999999998 -> -917.1819

Write-Host "Applied changes to all sheets"